$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Admins")

# Delete the rows for Elliot, Engineering, HumanResources and User (rows 4-7);
# only Admin and Artyom should remain as admin accounts.
$ws.Range("A4:C7").EntireRow.Delete()

# Swap the Admin (row 2) and Artyom (row 3) records so Artyom's record now comes
# first, ahead of Admin. Route through a scratch row so the copy preserves the
# underlying shared-string cell types instead of re-typing numeric-looking text.
$ws.Range("A2:C2").Copy($ws.Range("A100:C100"))
$ws.Range("A3:C3").Copy($ws.Range("A2:C2"))
$ws.Range("A100:C100").Copy($ws.Range("A3:C3"))
$ws.Range("A100:C100").Clear()
